$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns are treated as text so that
# numeric-looking strings (e.g. "35.323.23", "0.357") are preserved exactly
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '35.323.23'
$ws.Range("E2").Value = '  +0.53%  '

# Row 3
$ws.Range("D3").Value = '1.880.02'
$ws.Range("E3").Value = '  -1.21%  '

# Row 4
$ws.Range("E4").Value = '  -0.67%  '

# Row 5
$ws.Range("D5").Value = '247.33'
$ws.Range("E5").Value = '  -2.23%  '

# Row 6
$ws.Range("E6").Value = '  -2.35%  '

# Row 8
$ws.Range("D8").Value = '43.80'
$ws.Range("E8").Value = '  +5.11%  '

# Row 9
$ws.Range("D9").Value = '0.357'
$ws.Range("E9").Value = '  +1.11%  '

# Row 10
$ws.Range("D10").Value = '53.61'
$ws.Range("E10").Value = '  +2.29%  '

# Row 11
$ws.Range("D11").Value = '0.0741'
$ws.Range("E11").Value = '  -2.05%  '

# Row 12
$ws.Range("E12").Value = '  -0.22%  '

# Row 13
$ws.Range("D13").Value = '13.53'
$ws.Range("E13").Value = '  +2.22%  '

# Row 14
$ws.Range("D14").Value = '2.151.22'
$ws.Range("E14").Value = '  -1.33%  '

# Row 15
$ws.Range("D15").Value = '0.770'
$ws.Range("E15").Value = '  +4.96%  '

# Row 16
$ws.Range("D16").Value = '4.93'
$ws.Range("E16").Value = '  -1.34%  '

# Row 17
$ws.Range("D17").Value = '1.878.25'
$ws.Range("E17").Value = '  -1.50%  '

# Row 18
$ws.Range("D18").Value = '35.366.39'
$ws.Range("E18").Value = '  +0.61%  '

# Row 19
$ws.Range("D19").Value = '72.67'
$ws.Range("E19").Value = '  -1.37%  '

# Row 20
$ws.Range("E20").Value = '  -1.99%  '

# Row 21
$ws.Range("D21").Value = '244.14'
$ws.Range("E21").Value = '  +0.46%  '

# Row 22
$ws.Range("D22").Value = '12.88'
$ws.Range("E22").Value = '  -1.17%  '

# Row 23
$ws.Range("D23").Value = '5.00'
$ws.Range("E23").Value = '  -0.66%  '

# Row 24
$ws.Range("D24").Value = '2.66'
$ws.Range("E24").Value = '  +9.29%  '

# Row 26
$ws.Range("E26").Value = '  -2.25%  '

# Row 27
$ws.Range("D27").Value = '165.86'
$ws.Range("E27").Value = '  -1.75%  '

# Row 28
$ws.Range("D28").Value = '8.61'
$ws.Range("E28").Value = '  +0.26%  '

# Row 29
$ws.Range("D29").Value = '18.30'
$ws.Range("E29").Value = '  -1.13%  '

# Row 30
$ws.Range("E30").Value = '  -1.80%  '

# Row 31
$ws.Range("E31").Value = '  -2.19%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '4.30'
$ws.Range("E32").Value = '  -0.87%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.0592'
$ws.Range("E33").Value = '  -0.57%  '

# Row 34
$ws.Range("B34").Value = 'TrustWalletToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D34").Value = '1.63'
$ws.Range("E34").Value = '  +0.29%  '

# Row 35
$ws.Range("D35").Value = '4.17'
$ws.Range("E35").Value = '  -1.96%  '

# Row 36
$ws.Range("E36").Value = '  -0.74%  '

# Row 37
$ws.Range("D37").Value = '0.847'
$ws.Range("E37").Value = '  -0.27%  '

# Row 38
$ws.Range("D38").Value = '1.95'
$ws.Range("E38").Value = '  -3.24%  '

# Row 39
$ws.Range("D39").Value = '0.0734'
$ws.Range("E39").Value = '  +10.70%  '

# Row 40
$ws.Range("D40").Value = '17.53'
$ws.Range("E40").Value = '  +1.23%  '

# Row 41
$ws.Range("E41").Value = '  +1.66%  '

# Row 42
$ws.Range("D42").Value = '96.93'
$ws.Range("E42").Value = '  -0.62%  '

# Row 43
$ws.Range("E43").Value = '  -2.58%  '

# Row 44
$ws.Range("D44").Value = '1.312.56'
$ws.Range("E44").Value = '  +0.64%  '

# Row 45
$ws.Range("D45").Value = '2.39'
$ws.Range("E45").Value = '  -1.07%  '

# Row 46
$ws.Range("D46").Value = '0.0803'
$ws.Range("E46").Value = '  +6.68%  '

# Row 47
$ws.Range("E47").Value = '  -1.58%  '

# Row 48
$ws.Range("E48").Value = '  -0.93%  '

# Row 49
$ws.Range("D49").Value = '11.84'
$ws.Range("E49").Value = '  -0.82%  '

# Row 50
$ws.Range("D50").Value = '6.28'
$ws.Range("E50").Value = '  -4.49%  '

# Row 51
$ws.Range("E51").Value = '  -1.98%  '
